$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 872
$ws.Range("F3").Value = 1813
$ws.Range("F4").Value = 90
$ws.Range("G6").Value = "已停售"
$ws.Range("F7").Value = 1394
$ws.Range("F8").Value = 2126
$ws.Range("F9").Value = 987
$ws.Range("F11").Value = 2439
$ws.Range("F14").Value = 4094
$ws.Range("F16").Value = 382
$ws.Range("F17").Value = 3197
$ws.Range("F18").Value = 878
$ws.Range("F19").Value = 153
$ws.Range("F21").Value = 199
$ws.Range("F22").Value = 2150
$ws.Range("F23").Value = 1207
$ws.Range("F24").Value = 6
$ws.Range("F25").Value = 2007
$ws.Range("F26").Value = 419
$ws.Range("F29").Value = 8926
$ws.Range("F30").Value = 5813
$ws.Range("F31").Value = 362
$ws.Range("F33").Value = 775
$ws.Range("F34").Value = 32
$ws.Range("F35").Value = 788
$ws.Range("F36").Value = 3481
$ws.Range("F38").Value = 949
$ws.Range("F39").Value = 414
$ws.Range("F40").Value = 65
$ws.Range("F41").Value = 205
$ws.Range("F43").Value = 4689
$ws.Range("F44").Value = 5
$ws.Range("F45").Value = 900
$ws.Range("F47").Value = 423

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F17").Value = 3472

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 8458
$ws.Range("F3").Value = 380
$ws.Range("F4").Value = 1380

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 872
$ws.Range("F3").Value = 380
$ws.Range("F4").Value = 1380
$ws.Range("F5").Value = 90
$ws.Range("F7").Value = 1395
$ws.Range("F8").Value = 987
$ws.Range("F12").Value = 4094
$ws.Range("F14").Value = 382
$ws.Range("F15").Value = 3197
$ws.Range("F16").Value = 878
$ws.Range("F17").Value = 153
$ws.Range("F19").Value = 2150
$ws.Range("F22").Value = 1207
$ws.Range("F24").Value = 6
$ws.Range("F26").Value = 419
$ws.Range("F29").Value = 8926
$ws.Range("F30").Value = 3472
$ws.Range("F32").Value = 362
$ws.Range("F34").Value = 775
$ws.Range("F35").Value = 788
$ws.Range("F36").Value = 949
$ws.Range("F37").Value = 414
$ws.Range("F38").Value = 65
$ws.Range("F39").Value = 205
$ws.Range("F42").Value = 4689
$ws.Range("F43").Value = 5
$ws.Range("F44").Value = 900
$ws.Range("F46").Value = 424
